$d = $word.ActiveDocument

# Find the empty paragraph (paragraph 2: "Start of demonstration:" / <empty> / "Some value" / "End of demonstration.")
# and remove it entirely (including its paragraph mark), as AQL expressions
# that render to empty strings should not leave a stray empty line.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    # Paragraph text includes the trailing paragraph mark (CR); an empty
    # paragraph's Range.Text is just that CR character.
    if ($text -eq "`r") {
        $para.Range.Delete()
    }
}
